# Editing datasheet for privacy concerns
# Replace the real phone numbers in column B (rows 2-55) with a sequential,
# non-identifying placeholder series, then restore the view to the top of
# the sheet with D24 selected, matching the author's recorded selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startValue = 9876543210
$firstRow = 2
$lastRow = 55

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = $startValue + ($row - $firstRow)
}

# Restore the view: scroll back to the top-left and leave D24 selected,
# mirroring the saved sheetView/selection state in the edited workbook.
$ws.Range("D24").Select() | Out-Null

Write-Output "Replaced phone numbers in B$firstRow`:B$lastRow with sequential placeholders."
